# Apply similarity_scores.xlsx update: grow the similarity matrix from a
# 2x2 comparison (A1:C3) to a 4x4 comparison (A1:E5). One of the original
# document names changes, two new documents are added (columns D/E, rows
# 4/5) and every pairwise similarity score is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Update existing document label -----------------------------------------
# C1 (column header) and A3 (row header) both point at the same original
# shared string "CV_Canadien_anglais.pdf"; update both occurrences so the
# underlying text is replaced everywhere it is referenced.
$ws.Range("C1").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\model_linkedin.docx"
$ws.Range("A3").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\model_linkedin.docx"

# --- New column/row headers for the two newly-compared documents -----------
$ws.Range("D1").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\net2.pdf"
$ws.Range("E1").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\Ryan-Lattrel_App-Note.pdf"

$ws.Range("A4").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\net2.pdf"
$ws.Range("A5").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\Ryan-Lattrel_App-Note.pdf"

# Match the bold/bordered/centered header style already used on B1/C1/A2/A3.
$ws.Range("B1").Copy()
$ws.Range("D1:E1").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Updated similarity scores ----------------------------------------------
$ws.Range("B2").Value = 0.9999999403953554
$ws.Range("C2").Value = 0.4689624690312765
$ws.Range("D2").Value = 0.5073443268135294
$ws.Range("E2").Value = 0.4179190944842692

$ws.Range("B3").Value = 0.4689624690312765
$ws.Range("C3").Value = 1.000000119209291
$ws.Range("D3").Value = 0.4948627427259438
$ws.Range("E3").Value = 0.305423068565764

$ws.Range("B4").Value = 0.5073443268135294
$ws.Range("C4").Value = 0.4948627427259438
$ws.Range("D4").Value = 0.9999999403953582
$ws.Range("E4").Value = 0.5091269485357147

$ws.Range("B5").Value = 0.4179190944842692
$ws.Range("C5").Value = 0.305423068565764
$ws.Range("D5").Value = 0.5091269485357147
$ws.Range("E5").Value = 1.000000059604645
